$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '42.184.01'
$ws.Range("E2").Value = '  -2.03%  '

Set-TextValue "D3" '2.240.03'
$ws.Range("E3").Value = '  -2.92%  '

$ws.Range("E4").Value = '  +0.02%  '

Set-TextValue "D5" '245.99'
$ws.Range("E5").Value = '  -2.84%  '

Set-TextValue "D6" '0.631'
$ws.Range("E6").Value = '  -1.74%  '

Set-TextValue "D7" '76.15'
$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("E8").Value = '  +0.00%  '

Set-TextValue "D9" '0.619'
$ws.Range("E9").Value = '  -5.43%  '

Set-TextValue "D10" '41.50'
$ws.Range("E10").Value = '  +4.40%  '

Set-TextValue "D11" '0.0943'
$ws.Range("E11").Value = '  -4.82%  '

Set-TextValue "D12" '7.06'
$ws.Range("E12").Value = '  -8.70%  '

$ws.Range("E13").Value = '  -3.38%  '

Set-TextValue "D14" '2.571.11'
$ws.Range("E14").Value = '  -3.09%  '

Set-TextValue "D15" '14.62'
$ws.Range("E15").Value = '  -5.54%  '

Set-TextValue "D16" '0.853'
$ws.Range("E16").Value = '  -3.19%  '

Set-TextValue "D17" '2.235.02'
$ws.Range("E17").Value = '  -3.32%  '

Set-TextValue "D18" '41.964.51'
$ws.Range("E18").Value = '  -2.45%  '

$ws.Range("E19").Value = '  -3.63%  '

Set-TextValue "D20" '71.57'
$ws.Range("E20").Value = '  -2.10%  '

Set-TextValue "D21" '6.06'
$ws.Range("E21").Value = '  -3.86%  '

Set-TextValue "D22" '2.28'
$ws.Range("E22").Value = '  +0.97%  '

Set-TextValue "D23" '230.32'
$ws.Range("E23").Value = '  -3.55%  '

$ws.Range("E24").Value = '  +0.01%  '

Set-TextValue "D25" '3.69'
$ws.Range("E25").Value = '  -5.54%  '

Set-TextValue "D26" '11.19'
$ws.Range("E26").Value = '  -4.10%  '

$ws.Range("E27").Value = '  -6.20%  '

Set-TextValue "D28" '7.39'
$ws.Range("E28").Value = '  +15.39%  '

$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D29" '168.78'
$ws.Range("E29").Value = '  +0.46%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D30" '2.09'
$ws.Range("E30").Value = '  -4.81%  '

Set-TextValue "D31" '20.51'
$ws.Range("E31").Value = '  -3.62%  '

$ws.Range("E32").Value = '  -2.69%  '

Set-TextValue "D33" '32.14'
$ws.Range("E33").Value = '  +5.08%  '

$ws.Range("E34").Value = '  -7.41%  '

$ws.Range("E35").Value = '  -2.41%  '

Set-TextValue "D36" '4.44'
$ws.Range("E36").Value = '  -3.73%  '

Set-TextValue "D37" '4.93'
$ws.Range("E37").Value = '  +1.39%  '

Set-TextValue "D38" '0.0299'
$ws.Range("E38").Value = '  -4.76%  '

Set-TextValue "D39" '13.92'
$ws.Range("E39").Value = '  -0.16%  '

$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D40" '2.17'
$ws.Range("E40").Value = '  -8.18%  '

$ws.Range("B41").Value = 'THORChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue "D41" '5.81'
$ws.Range("E41").Value = '  -1.44%  '

Set-TextValue "D42" '112.47'
$ws.Range("E42").Value = '  +7.10%  '

Set-TextValue "D43" '0.201'
$ws.Range("E43").Value = '  -8.07%  '

Set-TextValue "D44" '60.45'
$ws.Range("E44").Value = '  -3.57%  '

Set-TextValue "D45" '8.65'
$ws.Range("E45").Value = '  -6.31%  '

$ws.Range("E46").Value = '  -4.46%  '

$ws.Range("E47").Value = '  -0.45%  '

Set-TextValue "D48" '1.13'
$ws.Range("E48").Value = '  -5.01%  '

$ws.Range("E49").Value = '  -2.73%  '

$ws.Range("B50").Value = 'FTXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue "D50" '4.18'
$ws.Range("E50").Value = '  -14.76%  '

Set-TextValue "D51" '2.25'
$ws.Range("E51").Value = '  -3.14%  '
